$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Strip trailing suffix codes from existing product codes in column B
$ws.Range("B1").Value = "2 32BUCA"
$ws.Range("B2").Value = "2 P1LING"
$ws.Range("B3").Value = "4 P1PAPP"
$ws.Range("B4").Value = "2 P1PAPP"

# Add a new row 5 with product info
$ws.Range("A5").Value = "bucatini`ncase`n1"
$ws.Range("B5").Value = "1 32BUCA"

# Avoid leaving a custom row height from the multi-line entry (keep default row sizing)
$ws.Rows(5).AutoFit()
